# Applies the "Update docx golden tests for style changes" edit:
#   1. Adds a new paragraph style "AbstractTitle" ("Abstract Title").
#   2. Changes the "Abstract" style's space-before from 300 to 100 (twips/20).
#   3. Gives the "ImportTok" character style a green, bold color.
#   4. Gives the "BuiltInTok" character style a green color.

$d = $word.ActiveDocument

# --- 1. New "Abstract Title" paragraph style -------------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)     # wdStyleTypeParagraph = 1
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1            # wdAlignParagraphCenter
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15          # 15pt -> w:before="300"

$abstractTitle.Font.Size = 10                            # -> w:sz="20"
$abstractTitle.Font.SizeBi = 10                           # -> w:szCs="20"
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060                        # BGR of 345A8A -> w:color="345A8A"

# --- 2. "Abstract" style spacing-before 300 -> 100 -------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5                 # 5pt -> w:before="100"

# --- 3. "ImportTok" character style gains bold + green color ---------------
$importTok = $d.Styles("ImportTok")
$importTok.Font.Color = 32768                              # BGR of 008000 -> w:color="008000"
$importTok.Font.Bold = $true

# --- 4. "BuiltInTok" character style gains green color ----------------------
$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 32768                              # BGR of 008000 -> w:color="008000"
